# Generate Report for Handoff
# Rotates the source file's GUID-based name (e2b00420-... -> 2175c337-...)
# across all three sheets, refreshes the associated xliff handoff file
# names / hashes, and bumps the handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "2175c337-6dc0-4918-aa9a-89e15c8d7752"

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-24 06:58:33"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.0009188c8570ccdc952443e13dc51e2934816f79.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 06:58:28"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.0009188c8570ccdc952443e13dc51e2934816f79.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 06:58:33"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
